$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9379826188087463
$ws.Range("B1").Value = 1.60292112827301
$ws.Range("C1").Value = 6.641695976257324
$ws.Range("D1").Value = 2.811967611312866
$ws.Range("E1").Value = 1.520861625671387
